$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 638.36
$ws.Range("J38").Value = 1610.25
$ws.Range("L38").Value = 4830.75
$ws.Range("N38").Value = -5574.75
$ws.Range("H43").Value = 2554.9375
$ws.Range("I43").Value = 2733.3333
$ws.Range("J43").Value = 2513.7693
$ws.Range("K43").Value = 2733.3333
$ws.Range("L43").Value = 2513.7693
$ws.Range("M43").Value = -2664.3333
$ws.Range("N43").Value = -2651.7693
$ws.Range("H64").Value = 2981.5454
$ws.Range("I64").Value = 2934
$ws.Range("J64").Value = 2999.375
$ws.Range("K64").Value = 2934
$ws.Range("L64").Value = 2999.375
$ws.Range("M64").Value = -2686
$ws.Range("N64").Value = -3495.375
$ws.Range("H67").Value = 2981.5454
$ws.Range("I67").Value = 2934
$ws.Range("J67").Value = 2999.375
$ws.Range("K67").Value = 2934
$ws.Range("L67").Value = 2999.375
$ws.Range("M67").Value = -2076
$ws.Range("N67").Value = -4715.375
$ws.Range("H69").Value = 66670240
$ws.Range("I69").Value = 3880
$ws.Range("J69").Value = 71432120
$ws.Range("K69").Value = 11640
$ws.Range("L69").Value = 214296360
$ws.Range("M69").Value = -10766
$ws.Range("N69").Value = -214298108
$ws.Range("H70").Value = 1623.5294
$ws.Range("I70").Value = 1161.5385
$ws.Range("J70").Value = 3125
$ws.Range("K70").Value = 3484.6155
$ws.Range("L70").Value = 9375
$ws.Range("M70").Value = -3214.6155
$ws.Range("N70").Value = -9915
$ws.Range("H72").Value = 66670240
$ws.Range("I72").Value = 3880
$ws.Range("J72").Value = 71432120
$ws.Range("K72").Value = 34920
$ws.Range("L72").Value = 642889080
$ws.Range("M72").Value = -30552
$ws.Range("N72").Value = -642897816
$ws.Range("H73").Value = 1623.5294
$ws.Range("I73").Value = 1161.5385
$ws.Range("J73").Value = 3125
$ws.Range("K73").Value = 3484.6155
$ws.Range("L73").Value = 9375
$ws.Range("M73").Value = -2548.6155
$ws.Range("N73").Value = -11247
$ws.Range("H76").Value = 3229.2942
$ws.Range("I76").Value = 3228.4285
$ws.Range("K76").Value = 3228.4285
$ws.Range("M76").Value = -2913.4285
$ws.Range("H79").Value = 3229.2942
$ws.Range("I79").Value = 3228.4285
$ws.Range("K79").Value = 3228.4285
$ws.Range("M79").Value = -2136.4285
$ws.Range("H86").Value = 166702180
$ws.Range("I86").Value = 166702180
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 166702180
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -166701057
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 166702180
$ws.Range("I89").Value = 166702180
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 833510900
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -833505284
$ws.Range("N89").ClearContents()
$ws.Range("H112").Value = 5231.8887
$ws.Range("I112").Value = 766.6667
$ws.Range("J112").Value = 5550.8335
$ws.Range("K112").Value = 2300.0001
$ws.Range("L112").Value = 16652.5005
$ws.Range("M112").Value = -1192.0001
$ws.Range("N112").Value = -18868.5005
$ws.Range("H128").Value = 45950
$ws.Range("J128").Value = 45950
$ws.Range("L128").Value = 45950
$ws.Range("N128").Value = -55910
$ws.Range("H129").Value = 1462.8518
$ws.Range("J129").Value = 1630.3043
$ws.Range("L129").Value = 4890.9129
$ws.Range("N129").Value = -14890.9129
$ws.Range("H132").Value = 3686.1875
$ws.Range("I132").Value = 3602.6897
$ws.Range("J132").Value = 4493.3335
$ws.Range("K132").Value = 10808.0691
$ws.Range("L132").Value = 13480.0005
$ws.Range("M132").Value = -8278.069100000001
$ws.Range("N132").Value = -18540.0005
$ws.Range("H138").Value = 4806.987
$ws.Range("I138").Value = 6812.5713
$ws.Range("J138").Value = 4606.4287
$ws.Range("K138").Value = 20437.7139
$ws.Range("L138").Value = 13819.2861
$ws.Range("M138").Value = -15297.7139
$ws.Range("N138").Value = -24099.2861
$ws.Range("H139").Value = 59615.383
$ws.Range("J139").Value = 59615.383
$ws.Range("L139").Value = 59615.383
$ws.Range("N139").Value = -69895.383
$ws.Range("H140").Value = 76441.25
$ws.Range("J140").Value = 76441.25
$ws.Range("L140").Value = 76441.25
$ws.Range("N140").Value = -86801.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 16
$ws.Range("H6").Value = 503
$ws.Range("J6").Value = 503
$ws.Range("L6").Value = 503
$ws.Range("N6").Value = -849
$ws.Range("H23").Value = 19999
$ws.Range("J23").Value = 19999
$ws.Range("L23").Value = 19999
$ws.Range("N23").Value = -20517
$ws.Range("H24").Value = 69924.336
$ws.Range("J24").Value = 69924.336
$ws.Range("L24").Value = 69924.336
$ws.Range("N24").Value = -70672.336
$ws.Range("H63").Value = 4650.7837
$ws.Range("J63").Value = 5136.6665
$ws.Range("L63").Value = 5136.6665
$ws.Range("N63").Value = -6508.6665
$ws.Range("H66").Value = 4650.7837
$ws.Range("J66").Value = 5136.6665
$ws.Range("L66").Value = 25683.3325
$ws.Range("N66").Value = -32547.3325
$ws.Range("H68").Value = 46250
$ws.Range("J68").Value = 46250
$ws.Range("L68").Value = 46250
$ws.Range("N68").Value = -47872
$ws.Range("H71").Value = 46250
$ws.Range("J71").Value = 46250
$ws.Range("L71").Value = 138750
$ws.Range("N71").Value = -146862
$ws.Range("H80").Value = 20999
$ws.Range("J80").Value = 20999
$ws.Range("L80").Value = 20999
$ws.Range("N80").Value = -22995
$ws.Range("H83").Value = 20999
$ws.Range("J83").Value = 20999
$ws.Range("L83").Value = 62997
$ws.Range("N83").Value = -72981
$ws.Range("H100").Value = 69924.336
$ws.Range("J100").Value = 69924.336
$ws.Range("L100").Value = 69924.336
$ws.Range("N100").Value = -72088.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 34500
$ws.Range("J69").Value = 34500
$ws.Range("L69").Value = 34500
$ws.Range("N69").Value = -36122
$ws.Range("H72").Value = 34500
$ws.Range("J72").Value = 34500
$ws.Range("L72").Value = 103500
$ws.Range("N72").Value = -111612
$ws.Range("H82").Value = 17107.572
$ws.Range("J82").Value = 21499
$ws.Range("L82").Value = 21499
$ws.Range("N82").Value = -22265
$ws.Range("H85").Value = 17107.572
$ws.Range("J85").Value = 21499
$ws.Range("L85").Value = 21499
$ws.Range("N85").Value = -24151
$ws.Range("H105").Value = 1762.0667
$ws.Range("I105").Value = 1916.625
$ws.Range("J105").Value = 1585.4286
$ws.Range("K105").Value = 1916.625
$ws.Range("L105").Value = 1585.4286
$ws.Range("M105").Value = -169.625
$ws.Range("N105").Value = -5079.4286

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9988.456
$ws.Range("I31").Value = 1163.5
$ws.Range("J31").Value = 13103.147
$ws.Range("K31").Value = 1163.5
$ws.Range("L31").Value = 13103.147
$ws.Range("M31").Value = -868.5
$ws.Range("N31").Value = -13693.147
$ws.Range("H34").Value = 9988.456
$ws.Range("I34").Value = 1163.5
$ws.Range("J34").Value = 13103.147
$ws.Range("K34").Value = 1163.5
$ws.Range("L34").Value = 13103.147
$ws.Range("M34").Value = -961.5
$ws.Range("N34").Value = -13507.147
$ws.Range("H122").Value = 2387.647
$ws.Range("J122").Value = 2506.4285
$ws.Range("L122").Value = 7519.2855
$ws.Range("N122").Value = -12419.2855
$ws.Range("H134").Value = 3474894
$ws.Range("I134").Value = 3790524
$ws.Range("J134").Value = 2964.8333
$ws.Range("K134").Value = 11371572
$ws.Range("L134").Value = 8894.499899999999
$ws.Range("M134").Value = -11369037
$ws.Range("N134").Value = -13964.4999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4574.6763
$ws.Range("I131").Value = 499
$ws.Range("J131").Value = 4829.4062
$ws.Range("K131").Value = 1497
$ws.Range("L131").Value = 14488.2186
$ws.Range("M131").Value = 3543
$ws.Range("N131").Value = -24568.2186
$ws.Range("H137").Value = 15929659
$ws.Range("I137").Value = 20609.834
$ws.Range("J137").Value = 22293278
$ws.Range("K137").Value = 61829.50199999999
$ws.Range("L137").Value = 66879834
$ws.Range("M137").Value = -56729.50199999999
$ws.Range("N137").Value = -66890034
$ws.Range("H140").Value = 1804.9706
$ws.Range("I140").Value = 1549.25
$ws.Range("J140").Value = 2032.2778
$ws.Range("K140").Value = 4647.75
$ws.Range("L140").Value = 6096.8334
$ws.Range("M140").Value = 532.25
$ws.Range("N140").Value = -16456.8334

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 514.25
$ws.Range("I43").Value = 514.25
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 514.25
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -363.25
$ws.Range("N43").ClearContents()
$ws.Range("H46").Value = 4307.615
$ws.Range("J46").Value = 4307.615
$ws.Range("L46").Value = 4307.615
$ws.Range("N46").Value = -4619.615
$ws.Range("H57").Value = 19000
$ws.Range("J57").Value = 19000
$ws.Range("L57").Value = 19000
$ws.Range("N57").Value = -20640
$ws.Range("H70").Value = 5082.5713
$ws.Range("I70").Value = 5057.5
$ws.Range("J70").Value = 5145.25
$ws.Range("K70").Value = 5057.5
$ws.Range("L70").Value = 5145.25
$ws.Range("M70").Value = -4787.5
$ws.Range("N70").Value = -5685.25
$ws.Range("H73").Value = 5082.5713
$ws.Range("I73").Value = 5057.5
$ws.Range("J73").Value = 5145.25
$ws.Range("K73").Value = 5057.5
$ws.Range("L73").Value = 5145.25
$ws.Range("M73").Value = -4121.5
$ws.Range("N73").Value = -7017.25
$ws.Range("H80").Value = 3134933.2
$ws.Range("I80").Value = 4502400
$ws.Range("K80").Value = 4502400
$ws.Range("M80").Value = -4501402
$ws.Range("H83").Value = 3134933.2
$ws.Range("I83").Value = 4502400
$ws.Range("K83").Value = 22512000
$ws.Range("M83").Value = -22507008
$ws.Range("H102").Value = 2153.8
$ws.Range("I102").Value = 1989.1428
$ws.Range("J102").Value = 2538
$ws.Range("K102").Value = 1989.1428
$ws.Range("L102").Value = 2538
$ws.Range("M102").Value = -367.1428000000001
$ws.Range("N102").Value = -5782
$ws.Range("H113").Value = 62971.11
$ws.Range("I113").Value = 86452.234
$ws.Range("J113").Value = 1920.2
$ws.Range("K113").Value = 86452.234
$ws.Range("L113").Value = 1920.2
$ws.Range("M113").Value = -84282.234
$ws.Range("N113").Value = -6260.2
$ws.Range("H122").Value = 1682.5
$ws.Range("I122").Value = 1483
$ws.Range("J122").Value = 2041.6
$ws.Range("K122").Value = 4449
$ws.Range("L122").Value = 6124.799999999999
$ws.Range("M122").Value = -1999
$ws.Range("N122").Value = -11024.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1888.5
$ws.Range("I68").Value = 1786.3636
$ws.Range("K68").Value = 1786.3636
$ws.Range("M68").Value = -1037.3636
$ws.Range("H71").Value = 1888.5
$ws.Range("I71").Value = 1786.3636
$ws.Range("K71").Value = 8931.817999999999
$ws.Range("M71").Value = -5187.817999999999
$ws.Range("H82").Value = 1488.619
$ws.Range("I82").Value = 1622.4445
$ws.Range("J82").Value = 1388.25
$ws.Range("K82").Value = 1622.4445
$ws.Range("L82").Value = 1388.25
$ws.Range("M82").Value = -1261.4445
$ws.Range("N82").Value = -2110.25
$ws.Range("H85").Value = 1488.619
$ws.Range("I85").Value = 1622.4445
$ws.Range("J85").Value = 1388.25
$ws.Range("K85").Value = 1622.4445
$ws.Range("L85").Value = 1388.25
$ws.Range("M85").Value = -374.4445000000001
$ws.Range("N85").Value = -3884.25
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H105").Value = 95000
$ws.Range("J105").Value = 95000
$ws.Range("L105").Value = 95000
$ws.Range("N105").Value = -101988
$ws.Range("H122").Value = 6098.2617
$ws.Range("I122").Value = 5556.9414
$ws.Range("J122").Value = 6466.36
$ws.Range("K122").Value = 16670.8242
$ws.Range("L122").Value = 19399.08
$ws.Range("M122").Value = -14220.8242
$ws.Range("N122").Value = -24299.08

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 21142.143
$ws.Range("J54").Value = 21142.143
$ws.Range("L54").Value = 21142.143
$ws.Range("N54").Value = -22182.143
$ws.Range("H62").Value = 326499.75
$ws.Range("I62").Value = 500999.5
$ws.Range("J62").Value = 152000
$ws.Range("K62").Value = 500999.5
$ws.Range("L62").Value = 152000
$ws.Range("M62").Value = -500375.5
$ws.Range("N62").Value = -153248
$ws.Range("H65").Value = 326499.75
$ws.Range("I65").Value = 500999.5
$ws.Range("J65").Value = 152000
$ws.Range("K65").Value = 2504997.5
$ws.Range("L65").Value = 760000
$ws.Range("M65").Value = -2501877.5
$ws.Range("N65").Value = -766240
$ws.Range("H81").Value = 3501.5667
$ws.Range("I81").Value = 3940.5715
$ws.Range("J81").Value = 3265.1794
$ws.Range("K81").Value = 7881.143
$ws.Range("L81").Value = 6530.3588
$ws.Range("M81").Value = -6820.143
$ws.Range("N81").Value = -8652.3588
$ws.Range("H84").Value = 3501.5667
$ws.Range("I84").Value = 3940.5715
$ws.Range("J84").Value = 3265.1794
$ws.Range("K84").Value = 39405.715
$ws.Range("L84").Value = 32651.794
$ws.Range("M84").Value = -34101.715
$ws.Range("N84").Value = -43259.794
$ws.Range("H122").Value = 2571.258
$ws.Range("I122").Value = 2100
$ws.Range("K122").Value = 6300
$ws.Range("M122").Value = -3850
